$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.10766077041626
$ws.Range("B1").Value = 3.742244958877563
$ws.Range("C1").Value = 4.376760959625244
$ws.Range("D1").Value = 1.88029682636261
$ws.Range("E1").Value = 1.315934300422668
